$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking price strings so Excel
# does not auto-convert them to numbers (source column is text).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply updated prices / volume percentages
$ws.Range('D2').Value = '61.095.40'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '3.403.16'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '571.48'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').Value = '142.39'
$ws.Range('E6').Value = '  -1.55%  '
$ws.Range('D7').Value = '3.404.42'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').Value = '7.54'
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').Value = '0.396'
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('D13').Value = '3.979.18'
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('D14').Value = '28.53'
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '3.396.20'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '61.131.65'
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').Value = '14.03'
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('D21').Value = '9.10'
$ws.Range('E21').Value = '  -4.15%  '
$ws.Range('D22').Value = '386.90'
$ws.Range('E22').Value = '  -2.00%  '
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').Value = '74.22'
$ws.Range('E24').Value = '  +1.52%  '
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('E26').Value = '  -3.12%  '
$ws.Range('D27').Value = '3.533.37'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = '7.44'
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('D31').Value = '8.02'
$ws.Range('E31').Value = '  -1.64%  '
$ws.Range('E32').Value = '  -0.77%  '
$ws.Range('E33').Value = '  -3.28%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  -0.99%  '
$ws.Range('D36').Value = '7.04'
$ws.Range('E36').Value = '  +0.46%  '
$ws.Range('D37').Value = '166.59'
$ws.Range('E37').Value = '  -0.51%  '
$ws.Range('D38').Value = '3.430.08'
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('E39').Value = '  -1.59%  '
$ws.Range('E40').Value = '  -3.67%  '
$ws.Range('D41').Value = '28.44'
$ws.Range('E41').Value = '  +3.20%  '
$ws.Range('D42').Value = '0.0781'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').Value = '0.782'
$ws.Range('E43').Value = '  -2.44%  '
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').Value = '42.21'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  -0.73%  '
$ws.Range('D47').Value = '1.68'
$ws.Range('E47').Value = '  -2.64%  '
$ws.Range('D48').Value = '1.14'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('D49').Value = '2.506.33'
$ws.Range('E49').Value = '  -3.62%  '
$ws.Range('D50').Value = '23.46'
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('D51').Value = '6.84'
$ws.Range('E51').Value = '  -0.98%  '
